$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column layout -------------------------------------------------
# Column A becomes the wide "Job Description" column, column B a
# narrower "Resume" column, and a new column C holds the numeric
# "Fitment-Rating".
$ws.Columns.Item(1).ColumnWidth = (255.28515625 - 5/6)
$ws.Columns.Item(2).ColumnWidth = (54.85546875 - 5/6)

# ---- Header row ------------------------------------------------------
$ws.Range("A1").Value = "Job Description"
$ws.Range("B1").Value = "Resume"
$ws.Range("C1").Value = "Fitment-Rating"

$jobDescSDE = "We are actively looking for a skilled and innovative Software Development Engineer (SDE) to join our team. As an SDE, you will be at the forefront of our software development efforts, responsible for crafting elegant solutions to complex problems using a range of programming languages and technologies, including python, C++ , java and C. Your role will involve collaborating with cross-functional teams, designing and coding software applications, and ensuring the scalability and robustness of our systems. We also require you to have leadership skills and communication skills. "
$jobDescFE = "We are looking for a Front End Engineer with at least 2 years of experience in developing scalable and user-friendly web applications. The successful candidate will be proficient in modern JavaScript frameworks and libraries, HTML, CSS, and responsive design principles. This role will contribute significantly to the creation and implementation of user interfaces for our web applications."

# ---- Data rows ---------------------------------------------------
# r, JobDescription, Resume text, Fitment rating, A-style ("plain"/"grey"/"wrap")
$rows = @(
  @(2,  $jobDescSDE, "I have prior experience in python java but no prior experience in c++ and not experienced in leadership.", 0.6, "plain"),
  @(3,  $jobDescSDE, "I have no prior experience in python java c++ but I have good communication skills.", 0.3, "grey"),
  @(4,  $jobDescSDE, "I have prior experience in python but have no prior experience with java, c++ ,c  leadership and communication", 0.5, "grey"),
  @(5,  $jobDescSDE, "I have prior experience in python java c++ and c and also have prior experience in leadership and communication.", 1, "grey"),
  @(6,  $jobDescSDE, "I have prior experience in python java and leadership but no prior experience in c++ and c", 0.7, "grey"),
  @(7,  $jobDescSDE, "I have prior experience in python java c and leadership but no prior experience in c++ and communicaton.", 0.8, "grey"),
  @(8,  $jobDescSDE, "I have prior experience in python java but no prior experience in c++ and c.", 0.65, "grey"),
  @(9,  $jobDescSDE, "I have no prior experience in python java leadership but I have experience in communication, c++ and c.", 0.45, "grey"),
  @(10, $jobDescSDE, "I have prior experience in python java c and leadership but no prior experience in communicaton.", 0.85, "grey"),
  @(11, $jobDescSDE, "I have prior experience in python java  communication c and leadership but no prior experience in c++.", 0.9, "grey"),
  @(12, $jobDescFE, "I have two years of experience with prior knowledge in HTML CSS and JavaScript frameworks.I also have experience with React and Angular.", 1, "wrap"),
  @(13, $jobDescFE, "I have two years of experience with prior knowledge in HTML CSS and JavaScript frameworks.But I do not have experience with React and Angular", 0.9, "wrap"),
  @(14, $jobDescFE, "I have two years of experience with prior knowledge in HTML CSS React and Angular. I do not have prior experience with JavaScript.", 0.8, "grey"),
  @(15, $jobDescFE, "I have one year of experience with prior knowledge in HTML CSS React and Angular and JavaScript.", 0.85, "grey"),
  @(16, $jobDescFE, "I have one year of experience with prior knowledge in HTML and JavaScript frameworks.I alo have experience with React.", 0.75, "grey"),
  @(17, $jobDescFE, "I have two years of experience with prior knowledge in CSS and JavaScript frameworks.I do not have experience with Angular and HTML", 0.7, "grey"),
  @(18, $jobDescFE, "I have one year of experience with prior knowledge in HTML CSS and JavaScript frameworks.I do not have experience with React and Angular.", 0.6, "grey"),
  @(19, $jobDescFE, "I have one year of experience with prior knowledge in HTML CSS .I also have experience with React.", 0.5, "grey"),
  @(20, $jobDescFE, "I have two years of experience with prior knowledge in HTML with no prior knowledge of CSS JavaScript frameworks,React and Angular.", 0.3, "grey")
)

foreach ($row in $rows) {
  $r = $row[0]
  $jobDesc = $row[1]
  $resumeText = $row[2]
  $rating = $row[3]
  $styleKind = $row[4]

  $aCell = $ws.Cells.Item($r, 1)
  $bCell = $ws.Cells.Item($r, 2)
  $cCell = $ws.Cells.Item($r, 3)

  $aCell.Value = $jobDesc
  $bCell.Value = $resumeText
  $cCell.Value = $rating
  $cCell.NumberFormat = "0%"

  if ($styleKind -eq "grey") {
    $aCell.Font.Color = 4473924
  } elseif ($styleKind -eq "wrap") {
    $aCell.WrapText = $true
    $ws.Rows.Item($r).RowHeight = 30.75
  }
}

# ---- Selection / view state ---------------------------------------
$ws.Range("B21").Select()
